$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells - copy formatting (bold, centered, border) from the
# existing header cell AC1 so AD1:AF1 match the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record for every data row (2-37): Wins=96, Losses=67, Ties=0
$lastRow = 37
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 96   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 67   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
